$d = $word.ActiveDocument

# Title paragraph (1st paragraph): merge "Questions:" " " "Solving" " " "exponential" " " "equations"
# into a single run "Questions: Solving exponential equations".
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range
$titleRange.Find.Execute("Questions: Solving exponential equations", $true, $false, $false, $false, $false,
                          $true, 1, $false, "Questions: Solving exponential equations", 2)

# Author paragraph (2nd paragraph): merge the split name runs into a single run
# "Zoë Gemmell, Isabella Lewis, Akshat Srivastava".
$authorPara = $d.Paragraphs.Item(2)
$authorRange = $authorPara.Range
$authorRange.Find.Execute("Zoë Gemmell, Isabella Lewis, Akshat Srivastava", $true, $false, $false, $false, $false,
                           $true, 1, $false, "Zoë Gemmell, Isabella Lewis, Akshat Srivastava", 2)

# Abstract paragraph (4th paragraph): merge the split word runs into a single run.
$abstractPara = $d.Paragraphs.Item(4)
$abstractRange = $abstractPara.Range
$abstractRange.Find.Execute("A selection of questions for the study guide on solving equations involving indices.",
                             $true, $false, $false, $false, $false,
                             $true, 1, $false,
                             "A selection of questions for the study guide on solving equations involving indices.", 2)
